# Fruta / hortaliza, semanal
#
# The weekly data refresh re-sorts the 12 price-observation rows (rows 2-13,
# columns D..T: Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad). Columns A..C
# and E..K are identical across every row, so only D..T are touched, but we
# move the whole row for correctness/robustness.
#
# Row r in the refreshed sheet receives the data that, before the refresh,
# lived in row $rowMap[r].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 13
$firstCol = 4   # column D
$lastCol = 20   # column T

$rowMap = @{
    2  = 7
    3  = 8
    4  = 10
    5  = 11
    6  = 2
    7  = 3
    8  = 9
    9  = 12
    10 = 13
    11 = 6
    12 = 4
    13 = 5
}

# Snapshot every row's values first so that writes to one destination row
# never corrupt the source data still needed for another row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
